$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 to I1:J1, then set header text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in I0 and IF data for rows 2-54
$iVals = @(7,8,9,5,9,8,8,8,6,6,9,9,9,5,7,7,7,8,4,7,9,9,9,9,9,9,9,9,9,10,9,9,9,8,9,9,9,8,9,9,9,9,8,8,9,9,4,8,3,6,4,4,2)
$jVals = @(7,8,9,5,10,8,8,8,6,6,9,9,9,5,7,7,7,8,5,7,9,9,9,9,9,9,10,9,9,10,9,9,9,9,9,9,9,9,9,9,9,9,8,8,9,9,4,8,3,6,4,4,2)
for ($r = 2; $r -le 54; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}

Write-Output "done"
